# Add a new worksheet "expedia_search" with a small search-form data table.
$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "expedia_search"

# Header row.
$ws2.Range("A1").Value = "Going"
$ws2.Range("B1").Value = "From_Date"
$ws2.Range("C1").Value = "End_Date"

# Data row. From_Date / End_Date are stored as text ("17"/"18"), so force
# text formatting before writing the values to avoid them becoming numbers.
$ws2.Range("A1:C2").NumberFormat = "@"
$ws2.Range("A2").Value = "Virginia Beach"
$ws2.Range("B2").Value = "17"
$ws2.Range("C2").Value = "18"

$ws2.Columns.Item(1).ColumnWidth = 11.8
$ws2.PageSetup.Orientation = 1
[void]$ws2.Range("B4").Select()

# Restore the original sheet's selection/active state.
[void]$wb.Worksheets.Item(1).Activate()
[void]$wb.Worksheets.Item(1).Range("B4").Select()
